$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Range("B2").Value = 0.4835600043714692
$ws.Range("C2").Value = 0.2318306736508191
$ws.Range("D2").Value = 0.02719723791139472
$ws.Range("E2").Value = 0.09970402780494325
$ws.Range("F2").Value = 3.785143198636646
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 2.124850447428841
$ws.Range("J2").Value = 0.2065888780661069
$ws.Range("K2").Value = 0.7140377030638092
$ws.Range("M2").Value = 0.2989512149090956
# Row 3 (A3 = 1)
$ws.Range("B3").Value = 0.4641018815268865
$ws.Range("C3").Value = 0.2245224549583327
$ws.Range("D3").Value = 0.02765857281840312
$ws.Range("E3").Value = 0.09945838426494902
$ws.Range("F3").Value = 3.735384344745228
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 2.096630956416888
$ws.Range("J3").Value = 0.2050118552615956
$ws.Range("K3").Value = 0.6874523524116682
$ws.Range("M3").Value = 0.2929052685063844
# Row 4 (A4 = 2)
$ws.Range("B4").Value = 0.4525560959199595
$ws.Range("C4").Value = 0.2201839731731923
$ws.Range("D4").Value = 0.02797041583493076
$ws.Range("E4").Value = 0.09935553210996595
$ws.Range("F4").Value = 3.705944817236244
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 2.079804044332178
$ws.Range("J4").Value = 0.20411618900485
$ws.Range("K4").Value = 0.6716794352369675
$ws.Range("M4").Value = 0.2894000296760666
# Row 5 (A5 = 3)
$ws.Range("B5").Value = 0.4479520672405783
$ws.Range("C5").Value = 0.2184533389793728
$ws.Range("D5").Value = 0.02810470511489882
$ws.Range("E5").Value = 0.09932569896604271
$ws.Range("F5").Value = 3.694227242069786
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 2.073072290808994
$ws.Range("J5").Value = 0.2037694556778575
$ws.Range("K5").Value = 0.6653901955912715
$ws.Range("M5").Value = 0.2880236786987638
# Row 6 (A6 = 4)
$ws.Range("B6").Value = 0.4471936702024948
$ws.Range("C6").Value = 0.2181682208900071
$ws.Range("D6").Value = 0.02812743995606226
$ws.Range("E6").Value = 0.09932147528647484
$ws.Range("F6").Value = 3.692298400389632
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 2.071962048272496
$ws.Range("J6").Value = 0.2037129834452145
$ws.Range("K6").Value = 0.6643542255510511
$ws.Range("M6").Value = 0.2877982814919555
# Row 7 (A7 = 5)
$ws.Range("B7").Value = 0.4524935955525962
$ws.Range("C7").Value = 0.2201604821868415
$ws.Range("D7").Value = 0.02797219767890802
$ws.Range("E7").Value = 0.09935508083622224
$ws.Range("F7").Value = 3.705785659844537
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 2.079712750517416
$ws.Range("J7").Value = 0.2041114389192487
$ws.Range("K7").Value = 0.6715940562963283
$ws.Range("M7").Value = 0.2893812569106906
# Row 8 (A8 = 6)
$ws.Range("B8").Value = 0.4767674811760969
$ws.Range("C8").Value = 0.2292798642111507
$ws.Range("D8").Value = 0.02735038820382663
$ws.Range("E8").Value = 0.0996093823866353
$ws.Range("F8").Value = 3.767755099136721
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 2.115016363759707
$ws.Range("J8").Value = 0.2060300412324878
$ws.Range("K8").Value = 0.7047567251963187
$ws.Range("M8").Value = 0.2968236136465308
# Row 9 (A9 = 7)
$ws.Range("B9").Value = 0.5275589310409714
$ws.Range("C9").Value = 0.2483486270587889
$ws.Range("D9").Value = 0.02635678120905283
$ws.Range("E9").Value = 0.100488155412453
$ws.Range("F9").Value = 3.898136612664985
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 2.188234985402829
$ws.Range("J9").Value = 0.210369382461927
$ws.Range("K9").Value = 0.7741668372296999
$ws.Range("M9").Value = 0.3130613617856284
# Row 10 (A10 = 8)
$ws.Range("B10").Value = 0.5668313854141047
$ws.Range("C10").Value = 0.2630901935058034
$ws.Range("D10").Value = 0.02576305880057106
$ws.Range("E10").Value = 0.101364984091024
$ws.Range("F10").Value = 3.99938547664928
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 2.244496169848105
$ws.Range("J10").Value = 0.2139107220162799
$ws.Range("K10").Value = 0.8278518844749101
$ws.Range("M10").Value = 0.3259962775448244
# Row 11 (A11 = 9)
$ws.Range("B11").Value = 0.585124954387112
$ws.Range("C11").Value = 0.2699574749775877
$ws.Range("D11").Value = 0.02552227101156745
$ws.Range("E11").Value = 0.101813989881844
$ws.Range("F11").Value = 4.046644292057124
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 2.270634803216552
$ws.Range("J11").Value = 0.2155988326672542
$ws.Range("K11").Value = 0.8528635374124178
$ws.Range("M11").Value = 0.332099801368571
# Row 12 (A12 = 10)
$ws.Range("B12").Value = 0.5921139744347101
$ws.Range("C12").Value = 0.2725812583593381
$ws.Range("D12").Value = 0.02543528149741192
$ws.Range("E12").Value = 0.1019912145989714
$ws.Range("F12").Value = 4.06471335746653
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 2.280611735955375
$ws.Range("J12").Value = 0.2162491875456993
$ws.Range("K12").Value = 0.8624199168862674
$ws.Range("M12").Value = 0.334442625899797
# Row 13 (A13 = 11)
$ws.Range("B13").Value = 0.5906060215104389
$ws.Range("C13").Value = 0.2720151420449781
$ws.Range("D13").Value = 0.02545383016536817
$ws.Range("E13").Value = 0.1019527262450346
$ws.Range("F13").Value = 4.060814149333993
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 2.278459513805473
$ws.Range("J13").Value = 0.2161086279405495
$ws.Range("K13").Value = 0.8603579961331889
$ws.Range("M13").Value = 0.3339366532589594
# Row 14 (A14 = 12)
$ws.Range("B14").Value = 0.5856987095644115
$ws.Range("C14").Value = 0.2701728677952246
$ws.Range("D14").Value = 0.02551503044882608
$ws.Range("E14").Value = 0.1018284261208962
$ws.Range("F14").Value = 4.048127370801126
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 2.271454029866092
$ws.Range("J14").Value = 0.215652115144735
$ws.Range("K14").Value = 0.8536480418441954
$ws.Range("M14").Value = 0.3322919145389562
# Row 15 (A15 = 13)
$ws.Range("B15").Value = 0.5827008685247677
$ws.Range("C15").Value = 0.2690474581331443
$ws.Range("D15").Value = 0.02555306262201285
$ws.Range("E15").Value = 0.1017532254551483
$ws.Range("F15").Value = 4.040378923100178
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 2.26717324288299
$ws.Range("J15").Value = 0.2153739343025336
$ws.Range("K15").Value = 0.849549076708314
$ws.Range("M15").Value = 0.3312885739297897
# Row 16 (A16 = 14)
$ws.Range("B16").Value = 0.5656444808412289
$ws.Range("C16").Value = 0.2626446543791587
$ws.Range("D16").Value = 0.02577938257428158
$ws.Range("E16").Value = 0.1013366481754616
$ws.Range("F16").Value = 3.996321203328307
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 2.242798953209416
$ws.Range("J16").Value = 0.2138019532113091
$ws.Range("K16").Value = 0.8262292010589363
$ws.Range("M16").Value = 0.3256018127380429
# Row 17 (A17 = 15)
$ws.Range("B17").Value = 0.5552906785334812
$ws.Range("C17").Value = 0.258758118840916
$ws.Range("D17").Value = 0.02592571136398547
$ws.Range("E17").Value = 0.1010939221877294
$ws.Range("F17").Value = 3.969601037238277
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 2.227986029632433
$ws.Range("J17").Value = 0.2128573586941584
$ws.Range("K17").Value = 0.8120744514123999
$ws.Range("M17").Value = 0.3221693587124719
# Row 18 (A18 = 16)
$ws.Range("B18").Value = 0.5493757586978063
$ws.Range("C18").Value = 0.2565378595208188
$ws.Range("D18").Value = 0.02601263490089778
$ws.Range("E18").Value = 0.1009590320358065
$ws.Range("F18").Value = 3.954345237214966
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 2.219517314212951
$ws.Range("J18").Value = 0.2123213127321151
$ws.Range("K18").Value = 0.8039885425694138
$ws.Range("M18").Value = 0.3202157490654116
# Row 19 (A19 = 17)
$ws.Range("B19").Value = 0.5473799909460695
$ws.Range("C19").Value = 0.2557887207707381
$ws.Range("D19").Value = 0.02604254025515829
$ws.Range("E19").Value = 0.1009141715127981
$ws.Range("F19").Value = 3.949199263870469
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 2.216658748730239
$ws.Range("J19").Value = 0.2121410633105398
$ws.Range("K19").Value = 0.8012603251427493
$ws.Range("M19").Value = 0.3195578363635789
# Row 20 (A20 = 18)
$ws.Range("B20").Value = 0.5563886858506351
$ws.Range("C20").Value = 0.2591702760877297
$ws.Range("D20").Value = 0.02590984901758731
$ws.Range("E20").Value = 0.101119272436982
$ws.Range("F20").Value = 3.972433752201198
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 2.229557579746981
$ws.Range("J20").Value = 0.2129571609460257
$ws.Range("K20").Value = 0.8135755005306464
$ws.Range("M20").Value = 0.3225326123559569
# Row 21 (A21 = 19)
$ws.Range("B21").Value = 0.5871384326596569
$ws.Range("C21").Value = 0.270713355407139
$ws.Range("D21").Value = 0.02549694085438503
$ws.Range("E21").Value = 0.101864740885464
$ws.Range("F21").Value = 4.051849080457458
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 2.273509568444325
$ws.Range("J21").Value = 0.2157859025637663
$ws.Range("K21").Value = 0.8556166096261393
$ws.Range("M21").Value = 0.3327741577681991
# Row 22 (A22 = 20)
$ws.Range("B22").Value = 0.6075944963419602
$ws.Range("C22").Value = 0.2783932432781171
$ws.Range("D22").Value = 0.02525150346697913
$ws.Range("E22").Value = 0.1023938854739619
$ws.Range("F22").Value = 4.104761278917834
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 2.3026942293925
$ws.Range("J22").Value = 0.2176993771941866
$ws.Range("K22").Value = 0.8835884842228552
$ws.Range("M22").Value = 0.3396515173967387
# Row 23 (A23 = 21)
$ws.Range("B23").Value = 0.5966438140930563
$ws.Range("C23").Value = 0.2742818807037395
$ws.Range("D23").Value = 0.02538027049113012
$ws.Range("E23").Value = 0.1021076379159673
$ws.Range("F23").Value = 4.076428459294021
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 2.287075646971473
$ws.Range("J23").Value = 0.2166721931129132
$ws.Range("K23").Value = 0.8686139677733138
$ws.Range("M23").Value = 0.3359641093924353
# Row 24 (A24 = 22)
$ws.Range("B24").Value = 0.555892159569197
$ws.Range("C24").Value = 0.2589838956934898
$ws.Range("D24").Value = 0.02591701166986837
$ws.Range("E24").Value = 0.1011077970850565
$ws.Range("F24").Value = 3.971152751688351
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 2.228846934220343
$ws.Range("J24").Value = 0.212912018502017
$ws.Range("K24").Value = 0.8128967147469552
$ws.Range("M24").Value = 0.3223683238544197
# Row 25 (A25 = 23)
$ws.Range("B25").Value = 0.5134758245853845
$ws.Range("C25").Value = 0.2430621716447092
$ws.Range("D25").Value = 0.02660154947621152
$ws.Range("E25").Value = 0.1002097878250297
$ws.Range("F25").Value = 3.861910656092732
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 2.167996707210733
$ws.Range("J25").Value = 0.2091335549174502
$ws.Range("K25").Value = 0.7549187692405894
$ws.Range("M25").Value = 0.3084923671637014
